# chore: update Sheets via scheduled runner
# Refreshes the market-price-derived columns (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ -> H:N) for the
# rows whose underlying market data changed, across the ALC/ARM/BSM/CRP/
# GSM/LTW/WVR leve-profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 3454.3076
$ws.Range("J88").Value = 4890.875
$ws.Range("L88").Value = 4890.875
$ws.Range("N88").Value = -5702.875
$ws.Range("H91").Value = 3454.3076
$ws.Range("J91").Value = 4890.875
$ws.Range("L91").Value = 4890.875
$ws.Range("N91").Value = -7698.875
$ws.Range("H113").Value = 3302.889
$ws.Range("I113").Value = 3062.2307
$ws.Range("K113").Value = 3062.2307
$ws.Range("M113").Value = 191.7692999999999
$ws.Range("H138").Value = 5176.3896
$ws.Range("I138").Value = 14237.363
$ws.Range("J138").Value = 3099.9167
$ws.Range("K138").Value = 42712.089
$ws.Range("L138").Value = 9299.750100000001
$ws.Range("M138").Value = -37572.089
$ws.Range("N138").Value = -19579.7501
$ws.Range("H141").Value = 2966.6667
$ws.Range("I141").Value = 1900
$ws.Range("K141").Value = 5700
$ws.Range("M141").Value = -520

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 51108.184
$ws.Range("I45").Value = 94398.17999999999
$ws.Range("K45").Value = 94398.17999999999
$ws.Range("M45").Value = -94021.17999999999
$ws.Range("H61").Value = 2343400.5
$ws.Range("I61").Value = 49894.91
$ws.Range("K61").Value = 49894.91
$ws.Range("M61").Value = -49682.91
$ws.Range("H74").Value = 546339
$ws.Range("I74").Value = 3239
$ws.Range("K74").Value = 3239
$ws.Range("M74").Value = -2365
$ws.Range("H77").Value = 546339
$ws.Range("I77").Value = 3239
$ws.Range("K77").Value = 16195
$ws.Range("M77").Value = -11827
$ws.Range("H97").Value = 36936.668
$ws.Range("I97").Value = 50399.5
$ws.Range("J97").Value = 10011
$ws.Range("K97").Value = 50399.5
$ws.Range("L97").Value = 10011
$ws.Range("M97").Value = -49903.5
$ws.Range("N97").Value = -11003
$ws.Range("H110").Value = 1594.8422
$ws.Range("I110").Value = 1475.25
$ws.Range("J110").Value = 2232.6667
$ws.Range("K110").Value = 1475.25
$ws.Range("L110").Value = 2232.6667
$ws.Range("M110").Value = 569.75
$ws.Range("N110").Value = -6322.6667
$ws.Range("H122").Value = 1112.8572
$ws.Range("I122").Value = 1149.1666
$ws.Range("J122").Value = 895
$ws.Range("K122").Value = 3447.4998
$ws.Range("L122").Value = 2685
$ws.Range("M122").Value = -997.4998000000001
$ws.Range("N122").Value = -7585
$ws.Range("H132").Value = 2293.673
$ws.Range("I132").Value = 2174.175
$ws.Range("K132").Value = 6522.525000000001
$ws.Range("M132").Value = -3992.525000000001
$ws.Range("H136").Value = 2343400.5
$ws.Range("I136").Value = 49894.91
$ws.Range("K136").Value = 149684.73
$ws.Range("M136").Value = -147134.73

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 9725.286
$ws.Range("I105").Value = 11165.6
$ws.Range("J105").Value = 6124.5
$ws.Range("K105").Value = 11165.6
$ws.Range("L105").Value = 6124.5
$ws.Range("M105").Value = -9418.6
$ws.Range("N105").Value = -9618.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 2130.4
$ws.Range("I5").Value = 2575.75
$ws.Range("J5").Value = 349
$ws.Range("K5").Value = 2575.75
$ws.Range("L5").Value = 349
$ws.Range("M5").Value = -2463.75
$ws.Range("N5").Value = -573
$ws.Range("H62").Value = 5499.4
$ws.Range("I62").Value = 5101
$ws.Range("J62").Value = 5765
$ws.Range("K62").Value = 5101
$ws.Range("L62").Value = 5765
$ws.Range("M62").Value = -4477
$ws.Range("N62").Value = -7013
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H65").Value = 5499.4
$ws.Range("I65").Value = 5101
$ws.Range("J65").Value = 5765
$ws.Range("K65").Value = 25505
$ws.Range("L65").Value = 28825
$ws.Range("M65").Value = -22385
$ws.Range("N65").Value = -35065
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H99").Value = 2753.5
$ws.Range("I99").Value = 2012
$ws.Range("K99").Value = 2012
$ws.Range("M99").Value = -514
$ws.Range("H112").Value = 42850
$ws.Range("J112").Value = 42850
$ws.Range("L112").Value = 42850
$ws.Range("N112").Value = -45804
$ws.Range("H126").Value = 2753.5
$ws.Range("I126").Value = 2012
$ws.Range("K126").Value = 6036
$ws.Range("M126").Value = -3566
$ws.Range("H132").Value = 58250.445
$ws.Range("I132").Value = 64906.75
$ws.Range("K132").Value = 194720.25
$ws.Range("M132").Value = -192190.25
$ws.Range("H134").Value = 1459.7949
$ws.Range("I134").Value = 1131.6129
$ws.Range("K134").Value = 3394.8387
$ws.Range("M134").Value = -859.8387000000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 89651.7
$ws.Range("I80").Value = 104029.58
$ws.Range("J80").Value = 41383.145
$ws.Range("K80").Value = 104029.58
$ws.Range("L80").Value = 41383.145
$ws.Range("M80").Value = -103031.58
$ws.Range("N80").Value = -43379.145
$ws.Range("H83").Value = 89651.7
$ws.Range("I83").Value = 104029.58
$ws.Range("J83").Value = 41383.145
$ws.Range("K83").Value = 520147.9
$ws.Range("L83").Value = 206915.725
$ws.Range("M83").Value = -515155.9
$ws.Range("N83").Value = -216899.725
$ws.Range("H126").Value = 4945.091
$ws.Range("I126").Value = 4865.6665
$ws.Range("J126").Value = 4974.875
$ws.Range("K126").Value = 14596.9995
$ws.Range("L126").Value = 14924.625
$ws.Range("M126").Value = -12126.9995
$ws.Range("N126").Value = -19864.625
$ws.Range("H132").Value = 829007.3
$ws.Range("I132").Value = 2271.6667
$ws.Range("K132").Value = 6815.000100000001
$ws.Range("M132").Value = -4285.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 91666.55499999999
$ws.Range("J20").Value = 99999
$ws.Range("L20").Value = 99999
$ws.Range("N20").Value = -100451
$ws.Range("H42").Value = 16723000
$ws.Range("J42").Value = 34000
$ws.Range("L42").Value = 34000
$ws.Range("N42").Value = -35126
$ws.Range("H46").Value = 8380.25
$ws.Range("I46").Value = 10215.538
$ws.Range("J46").Value = 6211.273
$ws.Range("K46").Value = 10215.538
$ws.Range("L46").Value = 6211.273
$ws.Range("M46").Value = -10027.538
$ws.Range("N46").Value = -6587.273
$ws.Range("H49").Value = 16723000
$ws.Range("J49").Value = 34000
$ws.Range("L49").Value = 34000
$ws.Range("N49").Value = -34294
$ws.Range("H122").Value = 2797.625
$ws.Range("I122").Value = 2322.5625
$ws.Range("K122").Value = 6967.6875
$ws.Range("M122").Value = -4517.6875
$ws.Range("H132").Value = 1814.0938
$ws.Range("I132").Value = 1111.8182
$ws.Range("K132").Value = 3335.4546
$ws.Range("M132").Value = -805.4546
$ws.Range("H136").Value = 1278.9762
$ws.Range("I136").Value = 1927.591
$ws.Range("J136").Value = 1048.8226
$ws.Range("K136").Value = 5782.772999999999
$ws.Range("L136").Value = 3146.4678
$ws.Range("M136").Value = -3232.772999999999
$ws.Range("N136").Value = -8246.4678
$ws.Range("H140").Value = 68474.75
$ws.Range("J140").Value = 68474.75
$ws.Range("L140").Value = 68474.75
$ws.Range("N140").Value = -78834.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 103674.2
$ws.Range("J81").Value = 202565.6
$ws.Range("L81").Value = 405131.2
$ws.Range("N81").Value = -407253.2
$ws.Range("H84").Value = 103674.2
$ws.Range("J84").Value = 202565.6
$ws.Range("L84").Value = 2025656
$ws.Range("N84").Value = -2036264
$ws.Range("H107").Value = 412.5238
$ws.Range("I107").Value = 401.3889
$ws.Range("K107").Value = 1204.1667
$ws.Range("M107").Value = 715.8333
$ws.Range("H122").Value = 2566.45
$ws.Range("I122").Value = 1796.875
$ws.Range("K122").Value = 5390.625
$ws.Range("M122").Value = -2940.625
